$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 22:04"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 5131140
$ws.Cells.Item(4, 3).Value = 35616
$ws.Cells.Item(4, 4).Value = 2624700
$ws.Cells.Item(4, 5).Value = 2341718
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 628
$ws.Cells.Item(4, 8).Value = 164722

# Row 8: Sudafrica
$ws.Cells.Item(8, 1).Value = "Sudafrica"
$ws.Cells.Item(8, 2).Value = 553188
$ws.Cells.Item(8, 3).Value = 7712
$ws.Cells.Item(8, 4).Value = 404568
$ws.Cells.Item(8, 5).Value = 138410
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 301
$ws.Cells.Item(8, 8).Value = 10210

# Row 27: Canada
$ws.Cells.Item(27, 1).Value = "Canada"
$ws.Cells.Item(27, 2).Value = 119197
$ws.Cells.Item(27, 3).Value = 212
$ws.Cells.Item(27, 4).Value = 103542
$ws.Cells.Item(27, 5).Value = 6679
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 6
$ws.Cells.Item(27, 8).Value = 8976

# Row 59: Argelia
$ws.Cells.Item(59, 1).Value = "Argelia"
$ws.Cells.Item(59, 2).Value = 34693
$ws.Cells.Item(59, 3).Value = 538
$ws.Cells.Item(59, 4).Value = 24083
$ws.Cells.Item(59, 5).Value = 9317
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 11
$ws.Cells.Item(59, 8).Value = 1293

# Row 68: Costa Rica
$ws.Cells.Item(68, 1).Value = "Costa Rica"
$ws.Cells.Item(68, 2).Value = 22802
$ws.Cells.Item(68, 3).Value = 721
$ws.Cells.Item(68, 4).Value = 7589
$ws.Cells.Item(68, 5).Value = 14985
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 10
$ws.Cells.Item(68, 8).Value = 228

# Row 69: Nepal
$ws.Cells.Item(69, 1).Value = "Nepal"
$ws.Cells.Item(69, 2).Value = 22592
$ws.Cells.Item(69, 3).Value = 378
$ws.Cells.Item(69, 4).Value = 16313
$ws.Cells.Item(69, 5).Value = 6206
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 3
$ws.Cells.Item(69, 8).Value = 73

# Row 70: Etiopia
$ws.Cells.Item(70, 1).Value = "Etiopia"
$ws.Cells.Item(70, 2).Value = 22253
$ws.Cells.Item(70, 3).Value = 801
$ws.Cells.Item(70, 4).Value = 9707
$ws.Cells.Item(70, 5).Value = 12156
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 10
$ws.Cells.Item(70, 8).Value = 390

# Row 76: Costa de Marfil
$ws.Cells.Item(76, 1).Value = "Costa de Marfil"
$ws.Cells.Item(76, 2).Value = 16620
$ws.Cells.Item(76, 3).Value = 96
$ws.Cells.Item(76, 4).Value = 12893
$ws.Cells.Item(76, 5).Value = 3623
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 104

# Row 97: Mauritania
$ws.Cells.Item(97, 1).Value = "Mauritania"
$ws.Cells.Item(97, 2).Value = 6510
$ws.Cells.Item(97, 3).Value = 12
$ws.Cells.Item(97, 4).Value = 5527
$ws.Cells.Item(97, 5).Value = 826
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 157

# Row 98: Paraguay
$ws.Cells.Item(98, 1).Value = "Paraguay"
$ws.Cells.Item(98, 2).Value = 6508
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = 5123
$ws.Cells.Item(98, 5).Value = 1316
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 69

# Row 124: Mali
$ws.Cells.Item(124, 1).Value = "Mali"
$ws.Cells.Item(124, 2).Value = 2565
$ws.Cells.Item(124, 3).Value = 4
$ws.Cells.Item(124, 4).Value = 1960
$ws.Cells.Item(124, 5).Value = 480
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 125

# Row 127: Mozambique
$ws.Cells.Item(127, 1).Value = "Mozambique"
$ws.Cells.Item(127, 2).Value = 2241
$ws.Cells.Item(127, 3).Value = 28
$ws.Cells.Item(127, 4).Value = 832
$ws.Cells.Item(127, 5).Value = 1393
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = 16

# Row 128: Lituania
$ws.Cells.Item(128, 1).Value = "Lituania"
$ws.Cells.Item(128, 2).Value = 2231
$ws.Cells.Item(128, 3).Value = 37
$ws.Cells.Item(128, 4).Value = 1668
$ws.Cells.Item(128, 5).Value = 482
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 81

# Row 135: Sierra Leona
$ws.Cells.Item(135, 1).Value = "Sierra Leona"
$ws.Cells.Item(135, 2).Value = 1895
$ws.Cells.Item(135, 3).Value = 8
$ws.Cells.Item(135, 4).Value = 1442
$ws.Cells.Item(135, 5).Value = 385
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 68

# Row 165: Aruba
$ws.Cells.Item(165, 1).Value = "Aruba"
$ws.Cells.Item(165, 2).Value = 509
$ws.Cells.Item(165, 3).Value = 113
$ws.Cells.Item(165, 4).Value = 114
$ws.Cells.Item(165, 5).Value = 392
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 3

# Row 166: Tanzania
$ws.Cells.Item(166, 1).Value = "Tanzania"
$ws.Cells.Item(166, 2).Value = 509
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 183
$ws.Cells.Item(166, 5).Value = 305
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 21

# Row 167: Taiwan
$ws.Cells.Item(167, 1).Value = "Taiwan"
$ws.Cells.Item(167, 2).Value = 479
$ws.Cells.Item(167, 3).Value = 2
$ws.Cells.Item(167, 4).Value = 443
$ws.Cells.Item(167, 5).Value = 29
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 7

# Row 168: Burundi
$ws.Cells.Item(168, 1).Value = "Burundi"
$ws.Cells.Item(168, 2).Value = 400
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 304
$ws.Cells.Item(168, 5).Value = 95
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 1

# Row 202: Timor Oriental
$ws.Cells.Item(202, 1).Value = "Timor Oriental"
$ws.Cells.Item(202, 2).Value = 25
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 24
$ws.Cells.Item(202, 5).Value = 1
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

# Row 203: Santa Lucia
$ws.Cells.Item(203, 1).Value = "Santa Lucia"
$ws.Cells.Item(203, 2).Value = 25
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 24
$ws.Cells.Item(203, 5).Value = 1
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0
